$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stahl_Treppe")

# The "Stundensatz (€) / P_Satz" row (row 9) is being removed from the
# calculation table; Excel shifts every row below it up by one.
$ws.Rows.Item(9).Delete()

# The final formula row (now row 10, previously row 11) no longer looks up
# the P_Satz variable - the hourly rate is hard-coded to 75 instead.
$ws.Range("E10").Value = "(math.ceil(H/0.18) * P_Stufe) + ((H * 1.8 * 2) * P_Wange * F_Faktor) + (L_Podest * B * P_Rost) + ((H/2.7) * 12 *75) + P_Mat"

# Make Stahl_Treppe the active sheet / tab and select the edited cell, so
# the workbook re-opens showing this sheet with E10 selected.
$ws.Activate()
$ws.Range("E10").Select()
